$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 255; this shifts existing rows 255-303 down to 256-304
$ws.Rows.Item(255).Insert()

# Populate the new row 255 with its values.
# Columns A,B,C,E,F,G,H,I,R are identical across all "Poroto granado" / "Primera" rows in this block,
# so copy them from row 256 (the row that used to be row 255 before the insert).
$ws.Cells.Item(255, 1).Value() = $ws.Cells.Item(256, 1).Value()
$ws.Cells.Item(255, 2).Value() = $ws.Cells.Item(256, 2).Value()
$ws.Cells.Item(255, 3).Value() = $ws.Cells.Item(256, 3).Value()
$ws.Cells.Item(255, 4).Value() = 44694
$ws.Cells.Item(255, 5).Value() = $ws.Cells.Item(256, 5).Value()
$ws.Cells.Item(255, 6).Value() = $ws.Cells.Item(256, 6).Value()
$ws.Cells.Item(255, 7).Value() = $ws.Cells.Item(256, 7).Value()
$ws.Cells.Item(255, 8).Value() = $ws.Cells.Item(256, 8).Value()
$ws.Cells.Item(255, 9).Value() = $ws.Cells.Item(256, 9).Value()
$ws.Cells.Item(255, 10).Value() = 52
$ws.Cells.Item(255, 11).Value() = 26000
$ws.Cells.Item(255, 12).Value() = 27000
$ws.Cells.Item(255, 13).Value() = 26500
$ws.Cells.Item(255, 14).Value() = "$/saco 25 kilos"
$ws.Cells.Item(255, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(255, 16).Value() = 1060
$ws.Cells.Item(255, 17).Value() = 25
$ws.Cells.Item(255, 18).Value() = $ws.Cells.Item(256, 18).Value()
